# Apply the data-cleaning edits described in the commit:
# - A handful of rows had their park_code corrected/anonymised and, where the
#   states/lat/long values no longer corresponded to a real park (placeholder
#   "xxxN" codes), those geo columns were cleared out entirely.
# - Other rows had their park_code corrected to the right park and the
#   states/lat/long updated to match that park's real location.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 127: Fort Caroline National Memorial -> code becomes placeholder xxx1,
# clear states/lat/long.
$ws.Range("C127").Value = "xxx1"
$ws.Range("E127:G127").ClearContents()

# Row 214: John D. Rockefeller -> code becomes placeholder xxx2, clear
# states/lat/long.
$ws.Range("C214").Value = "xxx2"
$ws.Range("E214:G214").ClearContents()

# Row 235: Lake Chelan National Recreation Area -> code becomes placeholder
# xxx3, clear states/lat/long.
$ws.Range("C235").Value = "xxx3"
$ws.Range("E235:G235").ClearContents()

# Row 260: Martin Luther King -> corrected park code "mlkm", and
# states/lat/long updated to the DC location.
$ws.Range("C260").Value = "mlkm"
$ws.Range("E260").Value = "DC"
$ws.Range("F260").Value = 38.8862276865
$ws.Range("G260").Value = -77.0442195534

# Row 298: Oregon Caves National Monument and Preserve -> corrected park code
# "orca", and states/lat/long updated to the OR location.
$ws.Range("C298").Value = "orca"
$ws.Range("E298").Value = "OR"
$ws.Range("F298").Value = 42.10319143
$ws.Range("G298").Value = -123.4018586

# Row 336: Ross Lake National Recreation Area -> code becomes placeholder
# xxx4, clear states/lat/long.
$ws.Range("C336").Value = "xxx4"
$ws.Range("E336:G336").ClearContents()

# Row 378: Timucuan Ecological and Historic Preserve -> corrected park code
# "nama", and states/lat/long updated to the DC location.
$ws.Range("C378").Value = "nama"
$ws.Range("E378").Value = "DC"
$ws.Range("F378").Value = 38.88162683
$ws.Range("G378").Value = -77.03586953

# Row 404: White House -> corrected park code "whho", and states/lat/long
# updated to the DC location.
$ws.Range("C404").Value = "whho"
$ws.Range("E404").Value = "DC"
$ws.Range("F404").Value = 38.89541886
$ws.Range("G404").Value = -77.03654147

# Row 412: World War I Memorial -> code becomes placeholder xxx6, clear
# states/lat/long.
$ws.Range("C412").Value = "xxx6"
$ws.Range("E412:G412").ClearContents()
